# Added New Mac-Address and Document Types
#
# Row 5's doctyp_code changes from "PSP" to "DOC001" (doccat_code stays "POI"),
# and a large block of new "doctyp_code"/"doccat_code" combinations is appended
# as rows 8-36 (rows 6 and 7 keep their original CRN/POR and COB/POB content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 5: PSP -> DOC001 (doccat_code "POI" unchanged)
$ws.Range("A5").Value = "DOC001"

# New rows 8-36: doctyp_code / doccat_code pairs to append below the
# existing data (rows 1-7). lang_code, is_active, cr_by and cr_dtimes are
# identical for every appended row.
$newRows = @(
    @(8,  "DOC001", "POI"),
    @(9,  "DOC002", "POI"),
    @(10, "DOC003", "POI"),
    @(11, "DOC004", "POI"),
    @(12, "DOC005", "POI"),
    @(13, "DOC006", "POI"),
    @(14, "DOC007", "POI"),
    @(15, "DOC008", "POI"),
    @(16, "DOC009", "POI"),
    @(17, "DOC010", "POI"),
    @(18, "DOC011", "POI"),
    @(19, "DOC012", "POI"),
    @(20, "DOC001", "POA"),
    @(21, "DOC013", "POA"),
    @(22, "DOC014", "POA"),
    @(23, "DOC015", "POA"),
    @(24, "DOC004", "POA"),
    @(25, "DOC005", "POA"),
    @(26, "DOC006", "POA"),
    @(27, "DOC016", "POA"),
    @(28, "DOC017", "POA"),
    @(29, "DOC018", "POA"),
    @(30, "DOC008", "POA"),
    @(31, "DOC024", "POR"),
    @(32, "DOC025", "POR"),
    @(33, "DOC026", "POR"),
    @(34, "DOC001", "POR"),
    @(35, "DOC027", "POR"),
    @(36, "DOC028", "POR")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $doctyp = $r[1]
    $doccat = $r[2]

    $ws.Cells.Item($rowNum, 1).Value = $doctyp      # doctyp_code
    $ws.Cells.Item($rowNum, 2).Value = $doccat       # doccat_code
    $ws.Cells.Item($rowNum, 3).Value = "ara"         # lang_code
    $ws.Cells.Item($rowNum, 4).Value = $true         # is_active
    $ws.Cells.Item($rowNum, 5).Value = "superadmin"  # cr_by
    $ws.Cells.Item($rowNum, 6).Value = "now()"       # cr_dtimes
}

# Reflect the new used-range selection state (active cell moves to G1,
# selection spans the rest of the columns to the right of the data).
$ws.Range("G1:XFD1048576").Select()
